# Update the two RAD test-run timestamps on the summary sheet:
# row 2 ("New Tax Return Amount Due" / Quarterly Estimated Tax run) and
# row 4 ("New Tax Return Amount Due" run) were re-executed; stamp their
# Date column with the new execution timestamps.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Sun Jan 14 16:33:59 EST 2024"
$ws.Range("B4").Value = "Sun Jan 14 16:34:13 EST 2024"
